$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("geometry")

# Insert a new column before AN, shifting AN:AQ -> AO:AR
$ws.Range("AN1:AN2").EntireColumn.Insert()

# Fill in the new column's header and value
$ws.Range("AO1").Copy()
$ws.Range("AN1").PasteSpecial(-4122)
$ws.Range("AN1").Value = "solidity"
$ws.Range("AN2").Value = "[1.42997704 1.70997375]"
